$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Risky / Revenue Growth / revenuegrowth / gt / (E=0.22)
$ws.Range("A5").Value = "Risky"
$ws.Range("B5").Value = "Revenue Growth"
$ws.Range("C5").Value = "revenuegrowth"
$ws.Range("D5").Value = "gt"
$ws.Range("E5").Value = 0.22

# Row 4: Defensive / Beta / beta / between / (F=1.1, G=0.5)
$ws.Range("A4").Value = "Defensive"
$ws.Range("B4").Value = "Beta"
$ws.Range("C4").Value = "beta"
$ws.Range("D4").Value = "between"
$ws.Range("F4").Value = 1.1
$ws.Range("G4").Value = 0.5

# Row 6: Risky / Beta / beta / gt / (E=1.1)
$ws.Range("A6").Value = "Risky"
$ws.Range("B6").Value = "Beta"
$ws.Range("C6").Value = "beta"
$ws.Range("D6").Value = "gt"
$ws.Range("E6").Value = 1.1

# Row 8: Standard / Market Cap / marketcap / gt / (E=2000000000)
$ws.Range("A8").Value = "Standard"
$ws.Range("B8").Value = "Market Cap"
$ws.Range("C8").Value = "marketcap"
$ws.Range("D8").Value = "gt"
$ws.Range("E8").Value = 2000000000

$ws.Range("G9").Select()
